$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with changed values.
# D-column values are forced to Text format first so Excel does not
# reinterpret decimal-looking strings (e.g. "144.82") as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.678.51"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.178.92"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.25"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.82"
$ws.Range("E6").Value = "  +4.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  +5.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.29"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("E10").Value = "  +4.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.429"
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.721.78"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.02"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000173"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.652.48"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.165.06"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.21"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.52"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.529"
$ws.Range("E23").Value = "  +4.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.30"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.92"
$ws.Range("E25").Value = "  +17.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.171"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0903"
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.91"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.34"
$ws.Range("E30").Value = "  +3.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.16"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.35"
$ws.Range("E32").Value = "  +3.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  +4.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.61"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("E36").Value = "  +4.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.745.71"
$ws.Range("E37").Value = "  +7.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.68"
$ws.Range("E40").Value = "  +2.55%  "
$ws.Range("E41").Value = "  +3.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.725"
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.45"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0290"
$ws.Range("E44").Value = "  +7.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.218.84"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.995"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.18"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.101"
$ws.Range("E48").Value = "  +6.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.50"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.772"
$ws.Range("E50").Value = "  +3.38%  "
$ws.Range("E51").Value = "  +0.00%  "

# Rows 38/39: EnergySwap and Hedera swapped places (with new price/volume figures)
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.39"
$ws.Range("E38").Value = "  -0.72%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0709"
$ws.Range("E39").Value = "  +5.76%  "
